$d = $word.ActiveDocument

$pairs = @(
    @{old="201×6=1206"; new="371×3=1113"},
    @{old="257×6=1542"; new="107×4=428"},
    @{old="411×7=2877"; new="751×8=6008"},
    @{old="191×2=382"; new="514×5=2570"},
    @{old="580×2=1160"; new="561×8=4488"},
    @{old="643×6=3858"; new="500×2=1000"},
    @{old="333×9=2997"; new="995×2=1990"},
    @{old="771×5=3855"; new="778×6=4668"},
    @{old="466×9=4194"; new="810×7=5670"},
    @{old="551×3=1653"; new="963×6=5778"},
    @{old="503×4=2012"; new="261×9=2349"},
    @{old="282×6=1692"; new="958×4=3832"},
    @{old="863×7=6041"; new="127×9=1143"},
    @{old="371×5=1855"; new="773×3=2319"},
    @{old="684×3=2052"; new="747×3=2241"},
    @{old="690×5=3450"; new="351×9=3159"},
    @{old="841×6=5046"; new="139×8=1112"},
    @{old="482×9=4338"; new="520×9=4680"},
    @{old="170×2=340"; new="894×6=5364"},
    @{old="192×6=1152"; new="873×7=6111"},
    @{old="524×2=1048"; new="566×3=1698"},
    @{old="119×9=1071"; new="536×2=1072"},
    @{old="507×4=2028"; new="465×7=3255"},
    @{old="316×4=1264"; new="773×8=6184"},
    @{old="316×9=2844"; new="376×7=2632"}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $pair.new, 2)
}

Write-Output "Done replacing $($pairs.Count) pairs"
